$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add a 7th sheet by copying the last existing sheet (keeps formatting /
#    column widths / styles identical to its siblings), placed after it.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)

# ---------------------------------------------------------------------------
# 2. For every sheet (the original 6 plus the freshly duplicated 7th):
#      - insert a new blank column at B (shifts old B..G to C..H)
#      - stamp the new column-B header value
#      - stamp the per-sheet "code" value that now lives in column C
#      - stamp the date value that now lives in the last column (H)
# ---------------------------------------------------------------------------
$codes = @(
    @("0120180600", "0120180605"),
    @("0120180606", "0120180611"),
    @("0120180612", "0120180617"),
    @("0120180618", "0120180623"),
    @("0120180624", "0120180629"),
    @("0120180630", "0120180635"),
    @("0120180636", "0120180641")
)

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    $ws.Columns("B:B").Insert()

    $ws.Range("B1").Value = "03"
    $ws.Range("C1").Value = $codes[$i - 1][0]
    $ws.Range("C2").Value = $codes[$i - 1][1]
    $ws.Range("H1").Value = "02/01/2018"
}

# ---------------------------------------------------------------------------
# 3. Rename every tab to its 1-based position ("1".."7").
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).Name = "$i"
}

# ---------------------------------------------------------------------------
# 4. Make the new 7th sheet the active tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(7).Activate()

# ---------------------------------------------------------------------------
# 5. Workbook / sheet VBA code names (best effort - harmless if unsupported).
# ---------------------------------------------------------------------------
$wb.CodeName = "EstaPasta_de_trabalho"
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $wb.Worksheets.Item($i).CodeName = "Planilha$i"
}
